{"js": "// Update the title date paragraph and the division-problem table cells\n// to match the new set of values from the commit.\n\nconst body = context.document.body;\n\n// --- 1. Title paragraph: \"2025-12-24 Wednesday\" -> \"2025-12-25 Thursday\" ---\nconst titlePara = body.paragraphs.getFirst();\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text === \"2025-12-24 Wednesday\") {\n  titlePara.insertText(\"2025-12-25 Thursday\", \"Replace\");\n}\n\n// --- 2. Table of division problems: 5 data rows (0,4,8,12,16) x 5 columns ---\nconst table = body.tables.getFirst();\n\n// Row -> column -> [oldText, newText]\nconst rowUpdates = [\n  { row: 0, cells: [\"37\u00f74=\", \"37\u00f79=\", \"72\u00f78=\", \"92\u00f75=\", \"63\u00f79=\"], newCells: [\"11\u00f74=\", \"89\u00f79=\", \"48\u00f73=\", \"27\u00f76=\", \"90\u00f76=\"] },\n  { row: 4, cells: [\"27\u00f72=\", \"56\u00f73=\", \"74\u00f76=\", \"10\u00f76=\", \"95\u00f72=\"], newCells: [\"14\u00f73=\", \"77\u00f77=\", \"45\u00f76=\", \"61\u00f79=\", \"32\u00f73=\"] },\n  { row: 8, cells: [\"40\u00f74=\", \"35\u00f76=\", \"73\u00f76=\", \"18\u00f78=\", \"79\u00f73=\"], newCells: [\"88\u00f78=\", \"65\u00f72=\", \"12\u00f73=\", \"51\u00f77=\", \"70\u00f77=\"] },\n  { row: 12, cells: [\"70\u00f75=\", \"54\u00f79=\", \"50\u00f72=\", \"54\u00f79=\", \"96\u00f79=\"], newCells: [\"21\u00f79=\", \"84\u00f79=\", \"50\u00f77=\", \"27\u00f72=\", \"47\u00f73=\"] },\n  { row: 16, cells: [\"93\u00f77=\", \"96\u00f73=\", \"88\u00f76=\", \"46\u00f75=\", \"56\u00f72=\"], newCells: [\"87\u00f77=\", \"11\u00f78=\", \"46\u00f76=\", \"89\u00f76=\", \"52\u00f73=\"] },\n];\n\nfor (const update of rowUpdates) {\n  for (let col = 0; col < update.newCells.length; col++) {\n    const cell = table.getCell(update.row, col);\n    const cellPara = cell.body.paragraphs.getFirst();\n    cellPara.insertText(update.newCells[col], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the title date paragraph and the division-problem table cells\n# to match the new set of values from the commit.\n\n$d = $word.ActiveDocument\n\n# --- 1. Title paragraph: \"2025-12-24 Wednesday\" -> \"2025-12-25 Thursday\" ---\n$titlePara = $d.Paragraphs.Item(1)\nif ($titlePara.Range.Text.TrimEnd(\"`r\") -eq \"2025-12-24 Wednesday\") {\n    $titleRange = $titlePara.Range\n    $titleRange.MoveEnd(1, -1) | Out-Null\n    $titleRange.Text = \"2025-12-25 Thursday\"\n}\n\n# --- 2. Table of division problems: 5 data rows (Word rows 1,5,9,13,17) x 5 columns ---\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @(\n    @{ Row = 1;  Old = @(\"37\u00f74=\", \"37\u00f79=\", \"72\u00f78=\", \"92\u00f75=\", \"63\u00f79=\"); New = @(\"11\u00f74=\", \"89\u00f79=\", \"48\u00f73=\", \"27\u00f76=\", \"90\u00f76=\") },\n    @{ Row = 5;  Old = @(\"27\u00f72=\", \"56\u00f73=\", \"74\u00f76=\", \"10\u00f76=\", \"95\u00f72=\"); New = @(\"14\u00f73=\", \"77\u00f77=\", \"45\u00f76=\", \"61\u00f79=\", \"32\u00f73=\") },\n    @{ Row = 9;  Old = @(\"40\u00f74=\", \"35\u00f76=\", \"73\u00f76=\", \"18\u00f78=\", \"79\u00f73=\"); New = @(\"88\u00f78=\", \"65\u00f72=\", \"12\u00f73=\", \"51\u00f77=\", \"70\u00f77=\") },\n    @{ Row = 13; Old = @(\"70\u00f75=\", \"54\u00f79=\", \"50\u00f72=\", \"54\u00f79=\", \"96\u00f79=\"); New = @(\"21\u00f79=\", \"84\u00f79=\", \"50\u00f77=\", \"27\u00f72=\", \"47\u00f73=\") },\n    @{ Row = 17; Old = @(\"93\u00f77=\", \"96\u00f73=\", \"88\u00f76=\", \"46\u00f75=\", \"56\u00f72=\"); New = @(\"87\u00f77=\", \"11\u00f78=\", \"46\u00f76=\", \"89\u00f76=\", \"52\u00f73=\") }\n)\n\nforeach ($update in $rowUpdates) {\n    $row = $update.Row\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $t.Cell($row, $col)\n        $cell.Range.Text = $update.New[$col - 1]\n    }\n}\n"}
